# Fruta / hortaliza, semanal
# Weekly update: a new price observation is inserted at row 21, pushing the
# existing observations (rows 21-42) down by one row (to rows 22-43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 21
$lastDataRow  = 42

# Shift existing rows down by one (process bottom-up so we never
# overwrite a row before it has been copied downward).
for ($r = $lastDataRow; $r -ge $firstDataRow; $r--) {
    $srcRange = $ws.Range("A" + $r + ":T" + $r)
    $values   = $srcRange.Value2
    $dstRange = $ws.Range("A" + ($r + 1) + ":T" + ($r + 1))
    $dstRange.Value2 = $values
}

# The newly created last row (old row 42's data, now at row 43) needs the
# same date number format as the other rows in column D.
$ws.Cells.Item($lastDataRow + 1, 4).NumberFormat = $ws.Cells.Item($firstDataRow, 4).NumberFormat

# Write the new observation into row 21 (only the fields that actually
# change; the rest of the row keeps the values it already had).
$ws.Cells.Item($firstDataRow, 4).Value2  = 44740   # Fecha
$ws.Cells.Item($firstDataRow, 14).Value2 = 34000   # Precio mínimo
$ws.Cells.Item($firstDataRow, 15).Value2 = 34000   # Precio máximo
$ws.Cells.Item($firstDataRow, 16).Value2 = 34000   # Precio promedio ponderado
$ws.Cells.Item($firstDataRow, 19).Value2 = 1889    # Precio $/Kg
